$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 7 (pushes old row 7 -> row 8, carrying its styles along)
$ws.Rows("7:7").Insert()

# New header row (mirrors row 2: UserName / Password / DashboardUrl)
$ws.Range("A7").Value = "UserName"
$ws.Range("B7").Value = "Password"
$ws.Range("C7").Value = "DashboardUrl"

# New data row (was row 7 before the insert): admin / admin / login URL + hyperlink
$ws.Range("C8").Value = "https://opensource-demo.orangehrmlive.com/web/index.php/auth/login"
$ws.Hyperlinks.Add($ws.Range("C8"), "https://opensource-demo.orangehrmlive.com/web/index.php/auth/login")

# Match the resulting selection left behind in the sheet view
[void]$ws.Range("C8").Select()
